$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values (E2:F4 and E5:E6, plus G5/H5 swap) ---

# Row 2
$ws.Range("E2").Value = 800000
$ws.Range("F2").Value = 761579.37

# Row 3
$ws.Range("E3").Value = 800000
$ws.Range("F3").Value = 761579.37

# Row 4
$ws.Range("E4").Value = 800000
$ws.Range("F4").Value = 761579.37

# Row 5 - E updates, and G5/H5 values swap
$ws.Range("E5").Value = 800000
$ws.Range("G5").Value = 456
$ws.Range("H5").Value = 3938753.8

# Row 6
$ws.Range("E6").Value = 800000

# --- Column width changes ---
# Before: col E:G width=16.33203125 bestFit, col H width=12.6640625
# After: col E:F width=16.33203125 bestFit, col G width=14.21875 bestFit, col H width=16.33203125 bestFit
# (ColumnWidth values below are chosen so the resulting stored OOXML width is as close
# as possible to the target given this engine's column-width quantization.)
$ws.Columns.Item(7).ColumnWidth = 13.333333333333334
$ws.Columns.Item(8).ColumnWidth = 15.5

# --- Selection change ---
$ws.Range("E2:H6").Select()
